# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update simple per-country numeric values (no row reordering) ---

# Alemania (row 8)
$ws.Range("B8").Value = 142283
$ws.Range("C8").Value = 886
$ws.Range("E8").Value = 52480
$ws.Range("G8").Value = 51
$ws.Range("H8").Value = 4403

# Moldavia (row 59)
$ws.Range("B59").Value = 2378
$ws.Range("C59").Value = 114
$ws.Range("E59").Value = 1930
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 57

# Uzbekistan (row 70)
$ws.Range("B70").Value = 1469
$ws.Range("C70").Value = 64
$ws.Range("E70").Value = 1271

# Mauricio (row 112)
$ws.Range("D112").Value = 180
$ws.Range("E112").Value = 136

# --- Re-sort Nepal / Mongolia / Sierra Leona block (rows 172-174) ---
# New order (by descending total cases): Nepal, Mongolia, Sierra Leona
# Nepal's own data also gets refreshed.

$ws.Range("A172").Value = "Nepal"
$ws.Range("B172").Value = 31
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 2
$ws.Range("E172").Value = 29
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

$ws.Range("A173").Value = "Mongolia"
$ws.Range("B173").Value = 31
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 5
$ws.Range("E173").Value = 26
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("A174").Value = "Sierra Leona"
$ws.Range("B174").Value = 30
$ws.Range("C174").Value = 4
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 30
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 16:52"
